$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 9680372
$ws.Range("J12").Value = 499.5
$ws.Range("L12").Value = 499.5
$ws.Range("N12").Value = -839.5
$ws.Range("H40").Value = 4103.3076
$ws.Range("I40").Value = 4309.75
$ws.Range("J40").Value = 4011.5557
$ws.Range("K40").Value = 4309.75
$ws.Range("L40").Value = 4011.5557
$ws.Range("M40").Value = -4134.75
$ws.Range("N40").Value = -4361.5557
$ws.Range("H42").Value = 1675.5385
$ws.Range("I42").Value = 298.5
$ws.Range("K42").Value = 895.5
$ws.Range("M42").Value = -665.5
$ws.Range("H43").Value = 15155849
$ws.Range("J43").Value = 5139.9375
$ws.Range("L43").Value = 5139.9375
$ws.Range("N43").Value = -5277.9375
$ws.Range("H48").Value = 400
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 400
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H92").Value = 2606873.2
$ws.Range("I92").Value = 2233103.2
$ws.Range("K92").Value = 2233103.2
$ws.Range("M92").Value = -2231855.2
$ws.Range("H98").Value = 1909.1111
$ws.Range("I98").Value = 1909.1111
$ws.Range("K98").Value = 1909.1111
$ws.Range("M98").Value = -411.1111000000001
$ws.Range("H122").Value = 1909.1111
$ws.Range("I122").Value = 1909.1111
$ws.Range("K122").Value = 5727.3333
$ws.Range("M122").Value = -3277.3333
$ws.Range("H125").Value = 455991.88
$ws.Range("I125").Value = 1421327.9
$ws.Range("K125").Value = 12791951.1
$ws.Range("M125").Value = -12789491.1
$ws.Range("H127").Value = 2455.25
$ws.Range("I127").Value = 1377.4286
$ws.Range("K127").Value = 4132.2858
$ws.Range("M127").Value = 827.7142000000003
$ws.Range("H132").Value = 2296.756
$ws.Range("I132").Value = 2248.5278
$ws.Range("K132").Value = 6745.5834
$ws.Range("M132").Value = -4215.5834
$ws.Range("H137").Value = 2129150.2
$ws.Range("I137").Value = 1031.875
$ws.Range("J137").Value = 3227534
$ws.Range("K137").Value = 3095.625
$ws.Range("L137").Value = 9682602
$ws.Range("M137").Value = -545.625
$ws.Range("N137").Value = -9687702

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 317349.53
$ws.Range("I32").Value = 317349.53
$ws.Range("K32").Value = 317349.53
$ws.Range("M32").Value = -317062.53
$ws.Range("H74").Value = 589117.2
$ws.Range("I74").Value = 2315.348
$ws.Range("J74").Value = 1488880.1
$ws.Range("K74").Value = 2315.348
$ws.Range("L74").Value = 1488880.1
$ws.Range("M74").Value = -1441.348
$ws.Range("N74").Value = -1490628.1
$ws.Range("H77").Value = 589117.2
$ws.Range("I77").Value = 2315.348
$ws.Range("J77").Value = 1488880.1
$ws.Range("K77").Value = 11576.74
$ws.Range("L77").Value = 7444400.5
$ws.Range("M77").Value = -7208.74
$ws.Range("N77").Value = -7453136.5
$ws.Range("H122").Value = 2210.1025
$ws.Range("I122").Value = 1521.6333
$ws.Range("J122").Value = 4505
$ws.Range("K122").Value = 4564.8999
$ws.Range("L122").Value = 13515
$ws.Range("M122").Value = -2114.8999
$ws.Range("N122").Value = -18415
$ws.Range("H140").Value = 107199.92
$ws.Range("J140").Value = 107199.92
$ws.Range("L140").Value = 107199.92
$ws.Range("N140").Value = -117559.92
$ws.Range("H141").Value = 103781.25
$ws.Range("J141").Value = 103781.25
$ws.Range("L141").Value = 103781.25
$ws.Range("N141").Value = -114141.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6641.7915
$ws.Range("I105").Value = 7675.8125
$ws.Range("K105").Value = 7675.8125
$ws.Range("M105").Value = -5928.8125
$ws.Range("H107").Value = 11575.407
$ws.Range("J107").Value = 4785.5
$ws.Range("L107").Value = 4785.5
$ws.Range("N107").Value = -8625.5
$ws.Range("H123").Value = 99888
$ws.Range("J123").Value = 99888
$ws.Range("L123").Value = 99888
$ws.Range("N123").Value = -109688
$ws.Range("H134").Value = 27274914
$ws.Range("I134").Value = 2005.4
$ws.Range("J134").Value = 300004000
$ws.Range("K134").Value = 6016.200000000001
$ws.Range("L134").Value = 900012000
$ws.Range("M134").Value = -3481.200000000001
$ws.Range("N134").Value = -900017070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2998.375
$ws.Range("J58").Value = 4050
$ws.Range("L58").Value = 4050
$ws.Range("N58").Value = -4456
$ws.Range("H86").Value = 145113.62
$ws.Range("I86").Value = 339571.66
$ws.Range("J86").Value = 28438.8
$ws.Range("K86").Value = 339571.66
$ws.Range("L86").Value = 28438.8
$ws.Range("M86").Value = -338448.66
$ws.Range("N86").Value = -30684.8
$ws.Range("H89").Value = 145113.62
$ws.Range("I89").Value = 339571.66
$ws.Range("J89").Value = 28438.8
$ws.Range("K89").Value = 1697858.3
$ws.Range("L89").Value = 142194
$ws.Range("M89").Value = -1692242.3
$ws.Range("N89").Value = -153426
$ws.Range("H107").Value = 1828.7059
$ws.Range("I107").Value = 1339.2
$ws.Range("J107").Value = 5500
$ws.Range("K107").Value = 1339.2
$ws.Range("L107").Value = 5500
$ws.Range("M107").Value = 580.8
$ws.Range("N107").Value = -9340
$ws.Range("H122").Value = 2621.2979
$ws.Range("I122").Value = 2822.5715
$ws.Range("K122").Value = 8467.7145
$ws.Range("M122").Value = -6017.7145
$ws.Range("H134").Value = 4150.3076
$ws.Range("I134").Value = 3833
$ws.Range("K134").Value = 11499
$ws.Range("M134").Value = -8964
$ws.Range("H136").Value = 2998.375
$ws.Range("J136").Value = 4050
$ws.Range("L136").Value = 12150
$ws.Range("N136").Value = -17250

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 612.9
$ws.Range("I10").Value = 354.83334
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 1064.50002
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -925.5000199999999
$ws.Range("N10").Value = -3278
$ws.Range("H33").Value = 73533.21000000001
$ws.Range("I33").Value = 997.5
$ws.Range("K33").Value = 5985
$ws.Range("M33").Value = -5702
$ws.Range("H38").Value = 701.9524
$ws.Range("J38").Value = 1289.8
$ws.Range("L38").Value = 3869.4
$ws.Range("N38").Value = -4563.4
$ws.Range("H39").Value = 5182.9
$ws.Range("I39").Value = 757.5
$ws.Range("K39").Value = 2272.5
$ws.Range("M39").Value = -1978.5
$ws.Range("H41").Value = 9900.666999999999
$ws.Range("I41").Value = 11180.8
$ws.Range("J41").Value = 3500
$ws.Range("K41").Value = 33542.39999999999
$ws.Range("L41").Value = 10500
$ws.Range("M41").Value = -33204.39999999999
$ws.Range("N41").Value = -11176
$ws.Range("H86").Value = 169150
$ws.Range("I86").Value = 2450
$ws.Range("J86").Value = 252500
$ws.Range("K86").Value = 7350
$ws.Range("L86").Value = 757500
$ws.Range("M86").Value = -6164
$ws.Range("N86").Value = -759872
$ws.Range("H89").Value = 169150
$ws.Range("I89").Value = 2450
$ws.Range("J89").Value = 252500
$ws.Range("K89").Value = 22050
$ws.Range("L89").Value = 2272500
$ws.Range("M89").Value = -16122
$ws.Range("N89").Value = -2284356
$ws.Range("H92").Value = 121.666664
$ws.Range("I92").Value = 121.666664
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 364.999992
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 883.000008
$ws.Range("N92").ClearContents()
$ws.Range("H137").Value = 6594.6523
$ws.Range("I137").Value = 5330.5
$ws.Range("J137").Value = 7040.8237
$ws.Range("K137").Value = 15991.5
$ws.Range("L137").Value = 21122.4711
$ws.Range("M137").Value = -10891.5
$ws.Range("N137").Value = -31322.4711

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5261.8096
$ws.Range("J2").Value = 9714.637000000001
$ws.Range("L2").Value = 9714.637000000001
$ws.Range("N2").Value = -9940.637000000001
$ws.Range("H70").Value = 7137.4546
$ws.Range("I70").Value = 7258.143
$ws.Range("K70").Value = 7258.143
$ws.Range("M70").Value = -6988.143
$ws.Range("H73").Value = 7137.4546
$ws.Range("I73").Value = 7258.143
$ws.Range("K73").Value = 7258.143
$ws.Range("M73").Value = -6322.143
$ws.Range("H102").Value = 15152431
$ws.Range("I102").Value = 16667604
$ws.Range("J102").Value = 699
$ws.Range("K102").Value = 16667604
$ws.Range("L102").Value = 699
$ws.Range("M102").Value = -16665982
$ws.Range("N102").Value = -3943
$ws.Range("H116").Value = 79443.55499999999
$ws.Range("J116").Value = 79443.55499999999
$ws.Range("L116").Value = 79443.55499999999
$ws.Range("N116").Value = -88621.55499999999
$ws.Range("H122").Value = 7792.6665
$ws.Range("I122").Value = 3945.7273
$ws.Range("J122").Value = 13837.857
$ws.Range("K122").Value = 11837.1819
$ws.Range("L122").Value = 41513.571
$ws.Range("M122").Value = -9387.1819
$ws.Range("N122").Value = -46413.571
$ws.Range("H126").Value = 2165.7693
$ws.Range("I126").Value = 2045.6
$ws.Range("K126").Value = 6136.799999999999
$ws.Range("M126").Value = -3666.799999999999
$ws.Range("H132").Value = 5106730.5
$ws.Range("I132").Value = 2434.4722
$ws.Range("J132").Value = 15315322
$ws.Range("K132").Value = 7303.4166
$ws.Range("L132").Value = 45945966
$ws.Range("M132").Value = -4773.4166
$ws.Range("N132").Value = -45951026

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 731.25
$ws.Range("J16").Value = 1088
$ws.Range("L16").Value = 1088
$ws.Range("N16").Value = -1428
$ws.Range("H22").Value = 1908.2941
$ws.Range("I22").Value = 1299
$ws.Range("J22").Value = 2162.1667
$ws.Range("K22").Value = 1299
$ws.Range("L22").Value = 2162.1667
$ws.Range("M22").Value = -1004
$ws.Range("N22").Value = -2752.1667
$ws.Range("H27").Value = 1908.2941
$ws.Range("I27").Value = 1299
$ws.Range("J27").Value = 2162.1667
$ws.Range("K27").Value = 1299
$ws.Range("L27").Value = 2162.1667
$ws.Range("M27").Value = -1192
$ws.Range("N27").Value = -2376.1667
$ws.Range("H40").Value = 4520.3687
$ws.Range("I40").Value = 4021.1428
$ws.Range("K40").Value = 4021.1428
$ws.Range("M40").Value = -3885.1428
$ws.Range("H122").Value = 4913.193
$ws.Range("I122").Value = 4015.125
$ws.Range("J122").Value = 6062.72
$ws.Range("K122").Value = 12045.375
$ws.Range("L122").Value = 18188.16
$ws.Range("M122").Value = -9595.375
$ws.Range("N122").Value = -23088.16

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8420.546
$ws.Range("I62").Value = 7865.75
$ws.Range("J62").Value = 8737.571
$ws.Range("K62").Value = 7865.75
$ws.Range("L62").Value = 8737.571
$ws.Range("M62").Value = -7241.75
$ws.Range("N62").Value = -9985.571
$ws.Range("H65").Value = 8420.546
$ws.Range("I65").Value = 7865.75
$ws.Range("J65").Value = 8737.571
$ws.Range("K65").Value = 39328.75
$ws.Range("L65").Value = 43687.855
$ws.Range("M65").Value = -36208.75
$ws.Range("N65").Value = -49927.855
$ws.Range("H81").Value = 4622.1113
$ws.Range("I81").Value = 2942.8572
$ws.Range("K81").Value = 5885.7144
$ws.Range("M81").Value = -4824.7144
$ws.Range("H84").Value = 4622.1113
$ws.Range("I84").Value = 2942.8572
$ws.Range("K84").Value = 29428.572
$ws.Range("M84").Value = -24124.572
$ws.Range("H100").Value = 779.8823
$ws.Range("I100").Value = 704.9167
$ws.Range("J100").Value = 959.8
$ws.Range("K100").Value = 1409.8334
$ws.Range("L100").Value = 1919.6
$ws.Range("M100").Value = -868.8334
$ws.Range("N100").Value = -3001.6
$ws.Range("H112").Value = 37075
$ws.Range("J112").Value = 37075
$ws.Range("L112").Value = 37075
$ws.Range("N112").Value = -40029
$ws.Range("H122").Value = 2293.8462
$ws.Range("I122").Value = 2305.85
$ws.Range("K122").Value = 6917.549999999999
$ws.Range("M122").Value = -4467.549999999999
$ws.Range("H132").Value = 65544.25
$ws.Range("I132").Value = 125863.625
$ws.Range("K132").Value = 377590.875
$ws.Range("M132").Value = -375060.875
$ws.Range("H136").Value = 32154.53
$ws.Range("I136").Value = 63611.688
$ws.Range("J136").Value = 4192.6113
$ws.Range("K136").Value = 190835.064
$ws.Range("L136").Value = 12577.8339
$ws.Range("M136").Value = -188285.064
$ws.Range("N136").Value = -17677.8339
